$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOM3228")

# --- "Ativação:" date 01/01/2012 -> 01/01/2023 ------------------------------
# B8/C8 and B15/C15 all share the same underlying string value. Plain
# `.Value = "01/01/2023"` gets re-parsed by Excel as a date literal (it would
# turn into the serial number 44927) because the text looks like a date.
# Route the write through a scratch cell holding a text *formula* (so the
# result is a genuine string, not a date) and PasteSpecial only the value
# into each destination cell - this keeps each cell's existing number format
# / style untouched and stores the text as a normal shared string, matching
# how the original cells were typed.
$scratch = $ws.Range("F1")
$scratch.Formula = "=""01/01/2023"""
$scratch.Copy()
$ws.Range("B8").PasteSpecial(-4163)
$ws.Range("C8").PasteSpecial(-4163)
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C15").PasteSpecial(-4163)
$scratch.ClearContents()

# --- "Objectives:" (row 11) - add the English translation ------------------
$ws.Range("B11").Value = "Provide knowledge about vacuum systems and production techniques and use of low temperatures."
$ws.Range("C11").Value = "Provide knowledge about vacuum systems and production techniques and use of low temperatures."

# --- "Short syllabus:" (row 14) - add the English translation --------------
$ws.Range("B14").Value = "Vacuum systems. Cryogenics and low temperature."
$ws.Range("C14").Value = "Vacuum systems. Cryogenics and low temperature."

# --- "Syllabus:" (row 16) - add the English translation --------------------
$ws.Range("B16").Value = "Theory of rarefied gases. Gas flow. Vacuum pumps. Quantitative description of the pumping of vacuum systems. Pressure gauges. Accessories: traps, shields, valves, etc. Adsorption, desorption and evaporation of molecules in vacuum. Leak detection .Sealing.Welding.Cleaning.cryogenics. Properties of cryogenic gases and liquids. Methods for obtaining low temperature. Liquefaction of gases. Temperature measurement. Cryogenic components. Calculation of heat transfer in cryostats and dewars."
$ws.Range("C16").Value = "Theory of rarefied gases. Gas flow. Vacuum pumps. Quantitative description of the pumping of vacuum systems. Pressure gauges. Accessories: traps, shields, valves, etc. Adsorption, desorption and evaporation of molecules in vacuum. Leak detection .Sealing.Welding.Cleaning.cryogenics. Properties of cryogenic gases and liquids. Methods for obtaining low temperature. Liquefaction of gases. Temperature measurement. Cryogenic components. Calculation of heat transfer in cryostats and dewars."

# --- Fix up cell styles for the newly-populated B/C cells ------------------
# These rows had no B/C cells before, so Excel falls back to the first
# matching <col> style entry (column A's bold style) instead of the correct
# per-column style. Re-apply the right look by copying formats from rows
# that already carry the correct column B / column C styles.
$ws.Range("B10").Copy()
$ws.Range("B11,B14,B16").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C11,C14,C16").PasteSpecial(-4122)
$excel.CutCopyMode = 0
